$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new value updates scraped from the latest cryptos.xlsx refresh.
# D-column cells whose new text would otherwise be parsed as a number
# (losing the original trailing-zero / fixed-precision formatting) are
# forced to Text format first so Excel stores the literal string.
$updates = @(
    @{ Cell = 'D2'; Value = '64.330.40'; AsText = $false }
    @{ Cell = 'E2'; Value = '  +0.56%  '; AsText = $false }
    @{ Cell = 'D3'; Value = '3.332.45'; AsText = $false }
    @{ Cell = 'E3'; Value = '  +0.27%  '; AsText = $false }
    @{ Cell = 'E4'; Value = '  -0.06%  '; AsText = $false }
    @{ Cell = 'D5'; Value = '553.60'; AsText = $true }
    @{ Cell = 'E5'; Value = '  +0.54%  '; AsText = $false }
    @{ Cell = 'D6'; Value = '173.37'; AsText = $true }
    @{ Cell = 'E6'; Value = '  +0.64%  '; AsText = $false }
    @{ Cell = 'D7'; Value = '0.621'; AsText = $true }
    @{ Cell = 'E7'; Value = '  +1.43%  '; AsText = $false }
    @{ Cell = 'D8'; Value = '1.00'; AsText = $true }
    @{ Cell = 'E8'; Value = '  -0.11%  '; AsText = $false }
    @{ Cell = 'D9'; Value = '3.323.29'; AsText = $false }
    @{ Cell = 'E9'; Value = '  +0.30%  '; AsText = $false }
    @{ Cell = 'E10'; Value = '  +5.85%  '; AsText = $false }
    @{ Cell = 'E11'; Value = '  +2.09%  '; AsText = $false }
    @{ Cell = 'D12'; Value = '53.49'; AsText = $true }
    @{ Cell = 'E12'; Value = '  +0.94%  '; AsText = $false }
    @{ Cell = 'E13'; Value = '  +2.84%  '; AsText = $false }
    @{ Cell = 'D14'; Value = '9.13'; AsText = $true }
    @{ Cell = 'E14'; Value = '  +1.65%  '; AsText = $false }
    @{ Cell = 'D15'; Value = '3.859.53'; AsText = $false }
    @{ Cell = 'E15'; Value = '  +0.02%  '; AsText = $false }
    @{ Cell = 'E16'; Value = '  +3.09%  '; AsText = $false }
    @{ Cell = 'D17'; Value = '18.13'; AsText = $true }
    @{ Cell = 'E17'; Value = '  -0.45%  '; AsText = $false }
    @{ Cell = 'D18'; Value = '3.327.03'; AsText = $false }
    @{ Cell = 'E18'; Value = '  +0.01%  '; AsText = $false }
    @{ Cell = 'D19'; Value = '64.277.99'; AsText = $false }
    @{ Cell = 'E19'; Value = '  +0.56%  '; AsText = $false }
    @{ Cell = 'D20'; Value = '11.74'; AsText = $true }
    @{ Cell = 'E20'; Value = '  +0.35%  '; AsText = $false }
    @{ Cell = 'D21'; Value = '0.987'; AsText = $true }
    @{ Cell = 'E21'; Value = '  +1.87%  '; AsText = $false }
    @{ Cell = 'D22'; Value = '453.72'; AsText = $true }
    @{ Cell = 'E22'; Value = '  +6.66%  '; AsText = $false }
    @{ Cell = 'D23'; Value = '5.08'; AsText = $true }
    @{ Cell = 'E23'; Value = '  +9.80%  '; AsText = $false }
    @{ Cell = 'E24'; Value = '  -0.36%  '; AsText = $false }
    @{ Cell = 'D25'; Value = '14.00'; AsText = $true }
    @{ Cell = 'E25'; Value = '  +5.45%  '; AsText = $false }
    @{ Cell = 'D26'; Value = '87.22'; AsText = $true }
    @{ Cell = 'E26'; Value = '  +3.92%  '; AsText = $false }
    @{ Cell = 'D27'; Value = '2.88'; AsText = $true }
    @{ Cell = 'E27'; Value = '  +2.81%  '; AsText = $false }
    @{ Cell = 'E28'; Value = '  +0.12%  '; AsText = $false }
    @{ Cell = 'D29'; Value = '30.99'; AsText = $true }
    @{ Cell = 'E29'; Value = '  +4.68%  '; AsText = $false }
    @{ Cell = 'D30'; Value = '8.59'; AsText = $true }
    @{ Cell = 'E30'; Value = '  +0.27%  '; AsText = $false }
    @{ Cell = 'D31'; Value = '6.53'; AsText = $true }
    @{ Cell = 'E31'; Value = '  -1.74%  '; AsText = $false }
    @{ Cell = 'D32'; Value = '11.42'; AsText = $true }
    @{ Cell = 'E32'; Value = '  +0.51%  '; AsText = $false }
    @{ Cell = 'D33'; Value = '61.75'; AsText = $true }
    @{ Cell = 'E33'; Value = '  +6.35%  '; AsText = $false }
    @{ Cell = 'D34'; Value = '566.79'; AsText = $true }
    @{ Cell = 'E34'; Value = '  -4.58%  '; AsText = $false }
    @{ Cell = 'E35'; Value = '  +0.53%  '; AsText = $false }
    @{ Cell = 'E36'; Value = '  +0.00%  '; AsText = $false }
    @{ Cell = 'E37'; Value = '  -1.10%  '; AsText = $false }
    @{ Cell = 'E38'; Value = '  +1.15%  '; AsText = $false }
    @{ Cell = 'D39'; Value = '35.38'; AsText = $true }
    @{ Cell = 'E39'; Value = '  +0.44%  '; AsText = $false }
    @{ Cell = 'D40'; Value = '0.366'; AsText = $true }
    @{ Cell = 'E40'; Value = '  +0.81%  '; AsText = $false }
    @{ Cell = 'E41'; Value = '  -2.35%  '; AsText = $false }
    @{ Cell = 'D42'; Value = '3.064.53'; AsText = $false }
    @{ Cell = 'E42'; Value = '  -0.66%  '; AsText = $false }
    @{ Cell = 'E43'; Value = '  +2.85%  '; AsText = $false }
    @{ Cell = 'E44'; Value = '  -0.95%  '; AsText = $false }
    @{ Cell = 'B45'; Value = 'Fetch.AI'; AsText = $false }
    @{ Cell = 'C45'; Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'; AsText = $false }
    @{ Cell = 'D45'; Value = '2.46'; AsText = $true }
    @{ Cell = 'E45'; Value = '  +1.15%  '; AsText = $false }
    @{ Cell = 'B46'; Value = 'Stellar'; AsText = $false }
    @{ Cell = 'C46'; Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'; AsText = $false }
    @{ Cell = 'D46'; Value = '0.134'; AsText = $true }
    @{ Cell = 'E46'; Value = '  +3.87%  '; AsText = $false }
    @{ Cell = 'B47'; Value = 'ApeXProtocol'; AsText = $false }
    @{ Cell = 'C47'; Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'; AsText = $false }
    @{ Cell = 'D47'; Value = '3.16'; AsText = $true }
    @{ Cell = 'E47'; Value = '  -1.44%  '; AsText = $false }
    @{ Cell = 'E48'; Value = '  -0.04%  '; AsText = $false }
    @{ Cell = 'D49'; Value = '140.81'; AsText = $true }
    @{ Cell = 'E49'; Value = '  +4.65%  '; AsText = $false }
    @{ Cell = 'E50'; Value = '  -3.33%  '; AsText = $false }
    @{ Cell = 'D51'; Value = '8.17'; AsText = $true }
    @{ Cell = 'E51'; Value = '  +0.99%  '; AsText = $false }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.AsText) {
        $rng.NumberFormat = "@"
    }
    $rng.Value = $u.Value
}
